$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relocate the "Status Types" reference column from E to G, keeping its
# formatting (header + the three status labels used for validation/legend).
$ws.Range("E1:E4").Cut($ws.Range("G1:G4"))
$ws.Range("E1:E4").Clear()

# Flesh out the two placeholder task rows now that the backend work has
# been scoped out - this is the "blank modal for register feature" work
# (Task 8 covers the companion Login backend).
$ws.Range("A9").Value = "Task 8: Login Backend"
$ws.Range("A10").Value = "Task 9: Register Backend"

# Restore the selection Excel had when the workbook was last saved.
$ws.Range("C14").Select()
